# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit: updated counters for
# several countries, the "Bahamas" / "Sri Lanka" row swap (Bahamas' updated
# totals now outrank Sri Lanka's, so it moves above it in the sorted table),
# and the "last updated" timestamp in the title cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 22:45"

# --- Bahamas / Sri Lanka swap (rows 128-129) -----------------------------
# Row 128 used to be Sri Lanka, row 129 used to be Bahamas. Bahamas' new
# totals (5191) now exceed Sri Lanka's (5170, unchanged), so Bahamas moves
# up to row 128 and Sri Lanka drops to row 129, carrying its old data with it.
$ws.Range("A128").Value = "Bahamas"
$ws.Range("B128").Value = 5191
$ws.Range("C128").Value = 28
$ws.Range("D128").Value = 3078
$ws.Range("E128").Value = 2004
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 109

$ws.Range("A129").Value = "Sri Lanka"
$ws.Range("B129").Value = 5170
$ws.Range("C129").Value = 132
$ws.Range("D129").Value = 3357
$ws.Range("E129").Value = 1800
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 13

# --- Row 4: Estados Unidos ------------------------------------------------
$ws.Range("B4").Value = 8138105
$ws.Range("C4").Value = 47852
$ws.Range("D4").Value = 5263370
$ws.Range("E4").Value = 2653033
$ws.Range("G4").Value = 829
$ws.Range("H4").Value = 221702

# --- Row 5: India ----------------------------------------------------------
$ws.Range("B5").Value = 7305060
$ws.Range("C5").Value = 67978
$ws.Range("D5").Value = 6378980
$ws.Range("E5").Value = 814786
$ws.Range("G5").Value = 677
$ws.Range("H5").Value = 111294

# --- Row 23: Peru ------------------------------------------------------
$ws.Range("B23").Value = 341512
$ws.Range("C23").Value = 5833
$ws.Range("E23").Value = 49845
$ws.Range("G23").Value = 27
$ws.Range("H23").Value = 9767

# --- Row 34 -----------------------------------------------------------------
$ws.Range("B34").Value = 149083
$ws.Range("C34").Value = 912
$ws.Range("E34").Value = 8685
$ws.Range("G34").Value = 29
$ws.Range("H34").Value = 12264

# --- Row 92 ------------------------------------------------------------
$ws.Range("B92").Value = 20217
$ws.Range("C92").Value = 34
$ws.Range("D92").Value = 19872
$ws.Range("E92").Value = 225

# --- Row 103 -----------------------------------------------------------
$ws.Range("B103").Value = 12069
$ws.Range("C103").Value = 69
$ws.Range("D103").Value = 10145
$ws.Range("E103").Value = 1794
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 130

# --- Row 109 -----------------------------------------------------------
$ws.Range("B109").Value = 10202
$ws.Range("C109").Value = 10
$ws.Range("D109").Value = 9892
$ws.Range("E109").Value = 241

# --- Row 118 -----------------------------------------------------------
$ws.Range("B118").Value = 7371
$ws.Range("C118").Value = 117
$ws.Range("D118").Value = 6270
$ws.Range("E118").Value = 1022
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 79

# --- Row 124 -----------------------------------------------------------
$ws.Range("B124").Value = 5715
$ws.Range("C124").Value = 19
$ws.Range("D124").Value = 5347
$ws.Range("E124").Value = 254

# --- Row 134 -----------------------------------------------------------
$ws.Range("B134").Value = 4940
$ws.Range("C134").Value = 32
$ws.Range("D134").Value = 4398
$ws.Range("E134").Value = 509
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 33
